$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time range text in cell B3 (17:05 - 17:10 -> 17:05 - 17:15)
$ws.Range("B3").Value = "17:05 - 17:15"

# Update the active selection to match the saved view state
$ws.Range("B4").Select()

$wb.Save()
